$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(
    -0.25566983192155135,
    -0.16045132858999622,
    -0.10970891086013701,
    -0.10146169206647571,
    -0.098032898476706265,
    -0.052714290666662222,
    -0.04210647134090717,
    -0.031953140488784726,
    -0.029664939539981727,
    -0.027422789285393989,
    -0.024406167329366291,
    -0.020878727214204584,
    -0.017177416190463646,
    -0.0090856242128962705,
    -0.0080555255634155998,
    -0.0060360759027298094,
    -0.0040037374720940377,
    -0.038923285815570807,
    -0.034859779180274941,
    -0.030420531537631135,
    -0.026305863645186278,
    -0.045714916945808426,
    -0.040500645690947579,
    -0.020099388704950094,
    -0.041865354966995838,
    -0.039347472000814676,
    -0.036827252025572577,
    -0.034719598041472288,
    -0.027654833264819878,
    0.032358114069944932,
    0.039449675831994213,
    0.01373359999162993
)

$colB = @(
    0.25511599810452168,
    0.15942382778356823,
    0.10946169202580514,
    0.10103289845312879,
    0.096579645607637765,
    0.052106471282517219,
    0.041953140428264923,
    0.031664939508263323,
    0.029422789251294823,
    0.027406167291085026,
    0.024378727173357539,
    0.020677416147119931,
    0.017085624151106593,
    0.0090555255287085856,
    0.0080360758636546237,
    0.0060037374323704817,
    0.0039999999522626339,
    0.038859779162837782,
    0.034420531516216268,
    0.030305863622949403,
    0.026146649338259387,
    0.045500645663331341,
    0.040099388616101628,
    0.019999999910014665,
    0.04184747197277261,
    0.039327251996898127,
    0.036719598011912602,
    0.034654833213654257,
    0.02764188566998893,
    -0.032449675881885742,
    -0.039514850691713121,
    -0.013750536784373324
)

for ($i = 0; $i -lt 32; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

# Column B width change (15.42578125 -> 14.7109375 stored width units)
$ws.Columns.Item(2).ColumnWidth = 13.8